$wb = $excel.ActiveWorkbook

# The "Repayment schedule" sheet is the 3rd sheet in the workbook.
$ws = $wb.Worksheets.Item(3)
$ws.Activate()

# Remember the width (in Excel's character-based ColumnWidth units) of the
# column to the left (M) so the newly-inserted column inherits the same
# width, mirroring Excel's default "insert column" behaviour.
$leftWidth = $ws.Columns("M").ColumnWidth

# Insert a new blank column before column N ("Late"), shifting the
# existing "Late", "heading"/Original and "Outstanding" columns one place
# to the right (N->O, O->P, P->Q).
$ws.Columns("N").Insert()

# Newly inserted column takes on the width of the column to its left.
$ws.Columns("N").ColumnWidth = $leftWidth

# Update the active selection left behind after the edit.
$ws.Range("R6").Select() | Out-Null
